# Update companies.xlsx: rename the two placeholder companies and append
# 13 new companies (rows 4-16), each with the same round-price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two existing placeholder company names.
$ws.Range("A2").Value = "MCWANN INFRA"
$ws.Range("A3").Value = "RUPANI CEMENT"

# New company names for rows 4-16.
$companyNames = @(
    "FIFADRA DEVELOPERS",
    "ROHIT ENTERPRISE",
    "NIMESH WATER PUMP",
    "AMITA INSURANCE",
    "CHARMACY",
    "9A LAB",
    "ZIUM LIFE SCIENCES",
    "GHELANI OIL CORP",
    "VINOD PETROLEUM",
    "VIVA GAS",
    "DOSHI PANTH PEDHI BANK",
    "HITEN CHIT FUND",
    "SUNITA SAHAKARI BANK"
)

$row = 4
foreach ($name in $companyNames) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 313
    $ws.Cells.Item($row, 3).Value = 1134
    $ws.Cells.Item($row, 4).Value = 3233
    $ws.Cells.Item($row, 5).Value = 234
    $ws.Cells.Item($row, 6).Value = 66
    $ws.Cells.Item($row, 7).Value = 456
    $row = $row + 1
}

# Update the selection to mirror the author's saved view state.
$ws.Range("G4:G16").Select()
